$wb = $excel.ActiveWorkbook

# Sheet1 holds the pending "names" queue; row 2 (id "qhnbcyxl") is being
# consumed/used, so remove it and let everything below shift up.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(2).Delete()

# "used" sheet gets a new trailing record documenting the consumed id.
$ws2 = $wb.Worksheets.Item("used")
$ws2.Range("A39").Value = "qhnbcyxl"
$ws2.Range("B39").Value = "ChatGPT Image 2026年1月21日 17_12_08.png"
$ws2.Range("C39").Value = "2026-01-21 17:12:54"
